$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45189 -> 45190) for every data row (rows 2 through 306).
for ($row = 2; $row -le 306; $row++) {
    $ws.Cells.Item($row, 3).Value = 45190
}
